# AFOLU general use tables — add derived summary statistics (Table 2.3)
# and fake emission-factor scratch data below the existing lookup table,
# and bold the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Bold the header row (A1 right-aligned, B1:F1 center-aligned — the
#    alignment was already set; we are only adding bold here).
# ---------------------------------------------------------------------
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1:F1").Font.Bold = $true

# ---------------------------------------------------------------------
# 2. New summary rows (12-17): average / min / max / min-ratio / max-ratio
#    per climate-region grouping, mirrored into column G as static values.
# ---------------------------------------------------------------------

# Row 12 - "dry"
$ws.Range("A12").Value = "dry"
$ws.Range("B12").Formula = "=AVERAGE(B3:F3,B5:F5,B7:F7)"
$ws.Range("C12").Formula = "=MIN(B3:F3,B5:F5,B7:F7)"
$ws.Range("D12").Formula = "=MAX(B3:F3,B5:F5,B7:F7)"
$ws.Range("E12").Formula = "=C12/B12"
$ws.Range("F12").Formula = "=D12/B12"
$ws.Range("G12").Value = 36.833333333333336

# Row 13 - "wet"
$ws.Range("A13").Value = "wet"
$ws.Range("B13").Formula = "=AVERAGE(B4:F4,B6:F6,B8:F8,B9:F9)"
$ws.Range("C13").Formula = "=MIN(B4:F4,B6:F6,B8:F8,B9:F9)"
$ws.Range("D13").Formula = "=MAX(B4:F4,B6:F6,B8:F8,B9:F9)"
$ws.Range("E13").Formula = "=C13/B13"
$ws.Range("F13").Formula = "=D13/B13"
$ws.Range("G13").Value = 75.411764705882348

# Row 14 - "tropical"
$ws.Range("A14").Value = "tropical"
$ws.Range("B14").Formula = "=AVERAGE(B7:F10)"
$ws.Range("C14").Formula = "=MIN(B7:F10)"
$ws.Range("D14").Formula = "=MAX(B7:F10)"
$ws.Range("E14").Formula = "=C14/B14"
$ws.Range("F14").Formula = "=D14/B14"
$ws.Range("G14").Value = 58.75
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = "KT SOC"

# Row 15 - "temperate_nutrient_rich"
$ws.Range("A15").Value = "temperate_nutrient_rich"
$ws.Range("B15").Formula = "=AVERAGE(B2:B6,E4,F2:F6)"
$ws.Range("C15").Formula = "=MIN(B3:B6,E4,F3:F6)"
$ws.Range("D15").Formula = "=MAX(B3:B6,E4,F3:F6)"
$ws.Range("E15").Formula = "=C15/B15"
$ws.Range("F15").Formula = "=D15/B15"
$ws.Range("G15").Value = 70.36363636363636
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = "KT N"

# Row 16 - "temperate_nutrient_poor"
$ws.Range("A16").Value = "temperate_nutrient_poor"
$ws.Range("B16").Formula = "=AVERAGE(C2:D6)"
$ws.Range("C16").Formula = "=MIN(C2:D6)"
$ws.Range("D16").Formula = "=MAX(C2:D6)"
$ws.Range("E16").Formula = "=C16/B16"
$ws.Range("F16").Formula = "=D16/B16"
$ws.Range("G16").Value = 41.444444444444443
$ws.Range("J16").Formula = "=10*0.0055"
$ws.Range("K16").Value = "KT N2O-N"

# Row 17 - "temperate"
$ws.Range("A17").Value = "temperate"
$ws.Range("B17").Formula = "=AVERAGE(B2:F6)"
$ws.Range("C17").Formula = "=MIN(B2:F6)"
$ws.Range("D17").Formula = "=MAX(B2:F6)"
$ws.Range("E17").Formula = "=C17/B17"
$ws.Range("F17").Formula = "=D17/B17"
$ws.Range("G17").Value = 60.19047619047619
$ws.Range("J17").Formula = "=J16*(11/7)*310"
$ws.Range("K17").Value = "KT N2O CO2E"
$ws.Range("L17").Value = "0.026 MT CO2E"

# ---------------------------------------------------------------------
# 3. Column A labels use the "right"-aligned style to match the rest of
#    the lookup table's first column (xlRight = -4152).
# ---------------------------------------------------------------------
$ws.Range("A12:A17").HorizontalAlignment = -4152

# ---------------------------------------------------------------------
# 4. Static duplicate scratch rows (19-24), mirroring rows 12-17's
#    E/F/G values as plain numbers (no formulas).
# ---------------------------------------------------------------------
$ws.Range("E19").Value = 0.51583710407239813
$ws.Range("F19").Value = 1.9004524886877827
$ws.Range("G19").Value = 36.833333333333336

$ws.Range("E20").Value = 0.45085803432137289
$ws.Range("F20").Value = 1.7238689547581905
$ws.Range("G20").Value = 75.411764705882348

$ws.Range("E21").Value = 0.52765957446808509
$ws.Range("F21").Value = 2.2127659574468086
$ws.Range("G21").Value = 58.75

$ws.Range("E22").Value = 0.28423772609819126
$ws.Range("F22").Value = 1.8475452196382429
$ws.Range("G22").Value = 70.36363636363636

$ws.Range("E23").Value = 0.2412868632707775
$ws.Range("F23").Value = 2.0509383378016088
$ws.Range("G23").Value = 41.444444444444443

$ws.Range("E24").Value = 0.16613924050632911
$ws.Range("F24").Value = 2.1598101265822787
$ws.Range("G24").Value = 60.19047619047619

# ---------------------------------------------------------------------
# 5. Selection moves to D9, matching the saved session state.
# ---------------------------------------------------------------------
$null = $ws.Range("D9").Select()
